$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 1025
$ws.Range("J45").Value = 1025
$ws.Range("L45").Value = 3075
$ws.Range("N45").Value = -3459
$ws.Range("H48").Value = 200
$ws.Range("J48").Value = 200
$ws.Range("L48").Value = 600
$ws.Range("N48").Value = -1184
$ws.Range("H56").Value = 200
$ws.Range("J56").Value = 200
$ws.Range("L56").Value = 600
$ws.Range("N56").Value = -1668
$ws.Range("H96").Value = 1504.1072
$ws.Range("I96").Value = 1619.5217
$ws.Range("J96").Value = 973.2
$ws.Range("K96").Value = 4858.5651
$ws.Range("L96").Value = 2919.6
$ws.Range("M96").Value = -3485.5651
$ws.Range("N96").Value = -5665.6
$ws.Range("H99").Value = 907.1667
$ws.Range("I99").Value = 953.2727
$ws.Range("J99").Value = 400
$ws.Range("K99").Value = 2859.8181
$ws.Range("L99").Value = 1200
$ws.Range("M99").Value = -1361.8181
$ws.Range("N99").Value = -4196
$ws.Range("H106").Value = 2781.3809
$ws.Range("I106").Value = 2993.2942
$ws.Range("J106").Value = 1880.75
$ws.Range("K106").Value = 2993.2942
$ws.Range("L106").Value = 1880.75
$ws.Range("M106").Value = -2362.2942
$ws.Range("N106").Value = -3142.75
$ws.Range("H112").Value = 4021.2727
$ws.Range("I112").Value = 2000
$ws.Range("J112").Value = 4223.4
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 12670.2
$ws.Range("M112").Value = -4892
$ws.Range("N112").Value = -14886.2
$ws.Range("H125").Value = 4657.625
$ws.Range("J125").Value = 860
$ws.Range("L125").Value = 7740
$ws.Range("N125").Value = -12660
$ws.Range("H132").Value = 14836.137
$ws.Range("I132").Value = 10780.429
$ws.Range("K132").Value = 32341.287
$ws.Range("M132").Value = -29811.287
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 3603.923
$ws.Range("I135").Value = 3954.6365
$ws.Range("J135").Value = 1675
$ws.Range("K135").Value = 35591.7285
$ws.Range("L135").Value = 15075
$ws.Range("M135").Value = -33056.7285
$ws.Range("N135").Value = -20145
$ws.Range("H138").Value = 2343.5757
$ws.Range("I138").Value = 1591.7646
$ws.Range("J138").Value = 2499.439
$ws.Range("K138").Value = 4775.293799999999
$ws.Range("L138").Value = 7498.316999999999
$ws.Range("M138").Value = 364.7062000000005
$ws.Range("N138").Value = -17778.317
$ws.Range("H140").Value = 224997.5
$ws.Range("J140").Value = 224997.5
$ws.Range("L140").Value = 224997.5
$ws.Range("N140").Value = -235357.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9257.105
$ws.Range("I2").Value = 7122.4287
$ws.Range("K2").Value = 7122.4287
$ws.Range("M2").Value = -7009.4287
$ws.Range("H32").Value = 5304.6514
$ws.Range("I32").Value = 3472.1875
$ws.Range("K32").Value = 3472.1875
$ws.Range("M32").Value = -3185.1875
$ws.Range("H61").Value = 124992.266
$ws.Range("I61").Value = 2015.8572
$ws.Range("K61").Value = 2015.8572
$ws.Range("M61").Value = -1803.8572
$ws.Range("H116").Value = 9257.105
$ws.Range("I116").Value = 7122.4287
$ws.Range("K116").Value = 7122.4287
$ws.Range("M116").Value = -4828.4287
$ws.Range("H136").Value = 124992.266
$ws.Range("I136").Value = 2015.8572
$ws.Range("K136").Value = 6047.571599999999
$ws.Range("M136").Value = -3497.571599999999
$ws.Range("H138").Value = 119198
$ws.Range("J138").Value = 119198
$ws.Range("L138").Value = 119198
$ws.Range("N138").Value = -129478
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9257.105
$ws.Range("I3").Value = 7122.4287
$ws.Range("K3").Value = 7122.4287
$ws.Range("M3").Value = -7008.4287
$ws.Range("H86").Value = 3197.7856
$ws.Range("I86").Value = 3052.2856
$ws.Range("K86").Value = 3052.2856
$ws.Range("M86").Value = -1929.2856
$ws.Range("H89").Value = 3197.7856
$ws.Range("I89").Value = 3052.2856
$ws.Range("K89").Value = 15261.428
$ws.Range("M89").Value = -9645.428
$ws.Range("H94").Value = 2421.8857
$ws.Range("I94").Value = 1497
$ws.Range("K94").Value = 1497
$ws.Range("M94").Value = -1046
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9192.4
$ws.Range("I16").Value = 8624.125
$ws.Range("K16").Value = 8624.125
$ws.Range("M16").Value = -8337.125
$ws.Range("H31").Value = 17373.242
$ws.Range("I31").Value = 7768.222
$ws.Range("K31").Value = 7768.222
$ws.Range("M31").Value = -7473.222
$ws.Range("H33").Value = 3685
$ws.Range("I33").Value = 2422
$ws.Range("K33").Value = 2422
$ws.Range("M33").Value = -2043
$ws.Range("H34").Value = 17373.242
$ws.Range("I34").Value = 7768.222
$ws.Range("K34").Value = 7768.222
$ws.Range("M34").Value = -7566.222
$ws.Range("H37").Value = 25000
$ws.Range("J37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("N37").Value = -25214
$ws.Range("H58").Value = 24776.709
$ws.Range("I58").Value = 8845
$ws.Range("K58").Value = 8845
$ws.Range("M58").Value = -8642
$ws.Range("H113").Value = 9192.4
$ws.Range("I113").Value = 8624.125
$ws.Range("K113").Value = 8624.125
$ws.Range("M113").Value = -6454.125
$ws.Range("H134").Value = 5784.1143
$ws.Range("I134").Value = 2092.923
$ws.Range("K134").Value = 6278.768999999999
$ws.Range("M134").Value = -3743.768999999999
$ws.Range("H136").Value = 24776.709
$ws.Range("I136").Value = 8845
$ws.Range("K136").Value = 26535
$ws.Range("M136").Value = -23985
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6212241.5
$ws.Range("I5").Value = 1141.1428
$ws.Range("J5").Value = 15873953
$ws.Range("K5").Value = 3423.4284
$ws.Range("L5").Value = 47621859
$ws.Range("M5").Value = -3311.4284
$ws.Range("N5").Value = -47622083
$ws.Range("H11").Value = 1192.909
$ws.Range("I11").Value = 1871.3334
$ws.Range("J11").Value = 378.8
$ws.Range("K11").Value = 5614.0002
$ws.Range("L11").Value = 1136.4
$ws.Range("M11").Value = -5474.0002
$ws.Range("N11").Value = -1416.4
$ws.Range("H23").Value = 682.25
$ws.Range("J23").Value = 1499
$ws.Range("L23").Value = 4497
$ws.Range("N23").Value = -4967
$ws.Range("H61").Value = 828.75
$ws.Range("I61").Value = 828.75
$ws.Range("K61").Value = 2486.25
$ws.Range("M61").Value = -2271.25
$ws.Range("H76").Value = 2007.6666
$ws.Range("I76").Value = 2007.6666
$ws.Range("K76").Value = 6022.9998
$ws.Range("M76").Value = -5639.9998
$ws.Range("H79").Value = 2007.6666
$ws.Range("I79").Value = 2007.6666
$ws.Range("K79").Value = 6022.9998
$ws.Range("M79").Value = -4696.9998
$ws.Range("H131").Value = 1353.52
$ws.Range("J131").Value = 1391.5714
$ws.Range("L131").Value = 4174.7142
$ws.Range("N131").Value = -14254.7142
$ws.Range("H135").Value = 6212241.5
$ws.Range("I135").Value = 1141.1428
$ws.Range("J135").Value = 15873953
$ws.Range("K135").Value = 10270.2852
$ws.Range("L135").Value = 142865577
$ws.Range("M135").Value = -7735.2852
$ws.Range("N135").Value = -142870647
$ws.Range("H139").Value = 11950.385
$ws.Range("I139").Value = 11950.385
$ws.Range("K139").Value = 35851.155
$ws.Range("M139").Value = -30711.155
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 19000
$ws.Range("J64").Value = 19000
$ws.Range("L64").Value = 19000
$ws.Range("N64").Value = -19450
$ws.Range("H67").Value = 19000
$ws.Range("J67").Value = 19000
$ws.Range("L67").Value = 19000
$ws.Range("N67").Value = -20560
$ws.Range("H93").Value = 6521.5356
$ws.Range("I93").Value = 5662
$ws.Range("K93").Value = 5662
$ws.Range("M93").Value = -4414
$ws.Range("H136").Value = 85446.66
$ws.Range("I136").Value = 145545.58
$ws.Range("J136").Value = 15331.25
$ws.Range("K136").Value = 436636.74
$ws.Range("L136").Value = 45993.75
$ws.Range("M136").Value = -434086.74
$ws.Range("N136").Value = -51093.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5544.2856
$ws.Range("I122").Value = 4564.273
$ws.Range("J122").Value = 7202.769
$ws.Range("K122").Value = 13692.819
$ws.Range("L122").Value = 21608.307
$ws.Range("M122").Value = -11242.819
$ws.Range("N122").Value = -26508.307
$ws.Range("H126").Value = 5189
$ws.Range("I126").Value = 5307.9
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 15923.7
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -13453.7
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 5473.683
$ws.Range("I132").Value = 2044.625
$ws.Range("J132").Value = 17665.889
$ws.Range("K132").Value = 6133.875
$ws.Range("L132").Value = 52997.667
$ws.Range("M132").Value = -3603.875
$ws.Range("N132").Value = -58057.667
$ws.Range("H136").Value = 11314.061
$ws.Range("I136").Value = 1063.3158
$ws.Range("J136").Value = 25225.785
$ws.Range("K136").Value = 3189.9474
$ws.Range("L136").Value = 75677.355
$ws.Range("M136").Value = -639.9474
$ws.Range("N136").Value = -80777.355
